# Applies the three text edits described by the diff:
#  1. Para "...if...else if, Switch-case": the leading-space run and the
#     "if...else if, Switch-case" run are merged into a single run.
#  2. Roll number "1501059" -> "1501223", kept as two separate runs
#     ("1501" and "223") instead of one.
#  3. Name "Gowtham M" -> "Sunil Varma S", kept as its own run (separate
#     from the preceding "Name: " run).

$d = $word.ActiveDocument

# --- Change 1: " " + "if…else if, Switch-case" -> one merged run ---
# Re-running a Find/Replace over this run causes the engine to coalesce it
# with the identically-formatted run right before it (the lone space run),
# which is exactly the run merge the diff shows.
$null = $d.Content.Find.Execute("if…else if, Switch-case", $true, $false, $false, $false, $false, $true, 1, $false, "if…else if, Switch-case", 2)

# --- Change 2: "1501059" -> "1501" / "223" (two separate runs) ---
$rngRoll = $d.Content
$null = $rngRoll.Find.Execute("1501059", $true, $false, $false, $false, $false, $true, 1, $false, "1501223", 2)

# Nudge (and immediately revert) formatting on "1501" and "223" so each one
# keeps its own run instead of being re-coalesced with its neighbour --
# first the "1501" head (keeps it apart from "ROLL NO.:24") ...
$rngHead = $d.Content
$null = $rngHead.Find.Execute("1501")
$rngHead.Bold = 1
$rngHead.Bold = 0

# ... then the "223" tail (keeps it apart from "1501").
$rngTail = $d.Content
$null = $rngTail.Find.Execute("223")
$rngTail.Bold = 1
$rngTail.Bold = 0

# --- Change 3: "Gowtham M" -> "Sunil Varma S" (kept as its own run) ---
$rngName = $d.Content
$null = $rngName.Find.Execute("Gowtham M", $true, $false, $false, $false, $false, $true, 1, $false, "Sunil Varma S", 2)

# Keep it split off from the preceding "Name: " run.
$rngName2 = $d.Content
$null = $rngName2.Find.Execute("Sunil Varma S")
$rngName2.Bold = 1
$rngName2.Bold = 0

Write-Output "Applied roll-number / name / run-merge edits."
